# Insert two new data rows at the top of the Ají price-record block
# (current rows 849-920 all shift down by two, becoming rows 851-922),
# then populate the two freshly inserted rows (849 and 850) with their
# own new record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("849:850").Insert()

# New row 849
$ws.Range("A849").Value = 6
$ws.Range("B849").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C849").Value = "Metropolitana"
$ws.Range("D849").Value = 44769
$ws.Range("E849").Value = 13
$ws.Range("F849").Value = 100112021
$ws.Range("G849").Value = "Ají"
$ws.Range("H849").Value = "Americana (o)"
$ws.Range("I849").Value = "Primera"
$ws.Range("J849").Value = 58
$ws.Range("K849").Value = 50000
$ws.Range("L849").Value = 55000
$ws.Range("M849").Value = 52241
$ws.Range("N849").Value = "$/caja 25 kilos"
$ws.Range("O849").Value = "Provincia de Limarí"
$ws.Range("P849").Value = 2090
$ws.Range("Q849").Value = 25
$ws.Range("R849").Value = "Hortaliza"

# New row 850
$ws.Range("A850").Value = 6
$ws.Range("B850").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C850").Value = "Metropolitana"
$ws.Range("D850").Value = 44769
$ws.Range("E850").Value = 13
$ws.Range("F850").Value = 100112021
$ws.Range("G850").Value = "Ají"
$ws.Range("H850").Value = "Americana (o)"
$ws.Range("I850").Value = "Segunda"
$ws.Range("J850").Value = 17
$ws.Range("K850").Value = 45000
$ws.Range("L850").Value = 45000
$ws.Range("M850").Value = 45000
$ws.Range("N850").Value = "$/caja 25 kilos"
$ws.Range("O850").Value = "Provincia de Limarí"
$ws.Range("P850").Value = 1800
$ws.Range("Q850").Value = 25
$ws.Range("R850").Value = "Hortaliza"
